# Pharma Society Report -- update membership counts and append a new
# survey response row (row 8), matching the authored diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated membership counts (column B) for existing societies ---
$ws.Range("B2").Value = 450
$ws.Range("B3").Value = 420
$ws.Range("B4").Value = 150
$ws.Range("B5").Value = 153
$ws.Range("B6").Value = 750

# Row 7, column B was stored as text "100"; normalize it to a real number.
$ws.Range("B7").Value = 100

# --- New row 8: another set of survey answers, appended below row 7 ---
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "500"
$ws.Range("B8").Style = "Normal"

$ws.Range("C8").Value = "Yes, community sites provide a platform for users to come together and interact."
$ws.Range("D8").Value = "Yes, research, Research can provide evidence and data to policymakers that can influence the development of state or local policies."
$ws.Range("E8").Value = "Yes, it allows for interaction and open communication with leaders."
$ws.Range("F8").Value = "Yes, the platform can target specific patient populations and help streamline the recruitment process."
$ws.Range("G8").Value = "Yes, it allows for direct interaction and collaboration with payors."
$ws.Range("H8").Value = "No, the board does not include area experts. Justification: lack of specific knowledge or expertise in the relevant field."
$ws.Range("I8").Value = "Yes, `nThere are many pharmaceutical companies and research organizations involved in therapeutic research collaborations to advance medical knowledge and improve patient outcomes."
$ws.Range("J8").Value = "No, there are no top therapeutic area experts on the board. This can be seen from the lack of expert professionals in the field serving on the board."
$ws.Range("K8").Value = "Asia."
